# Update the cached "datetimeFigureOut" footer-date text from 3/19/19 to
# 4/5/19 across the slide master, every slide layout, and the notes
# master (PowerPoint recalculates/recaches these on every save, which is
# what the original commit's diff captured).
#
# Also fix a typo in the "Storage" interface shape on slide 1:
# FoodDIaryStorage -> FoodDiaryStorage.

$p = $ppt.ActivePresentation

function Update-DateShapeText {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "3/19/19") {
                $tr.Text = "4/5/19"
                return $true
            }
        }
    }
    return $false
}

# Slide master date placeholder.
Update-DateShapeText $p.SlideMaster.Shapes | Out-Null

# Every slide layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DateShapeText $layout.Shapes | Out-Null
}

# Notes master date placeholder: its Shapes collection does not accept
# text writes in this host, but the HeadersFooters.DateAndTime proxy
# does, so use that instead.
$p.NotesMaster.HeadersFooters.DateAndTime.Text = "4/5/19"

# Fix "FoodDIaryStorage" -> "FoodDiaryStorage" on slide 1, touching only
# the misspelled run so the rest of the shape's text/formatting is left
# untouched.
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $fullRange = $shape.TextFrame.TextRange
        $fullText = $fullRange.Text
        $pos = $fullText.IndexOf("FoodDIaryStorage")
        if ($pos -ge 0) {
            $target = $fullRange.Characters($pos + 1, "FoodDIaryStorage".Length)
            $target.Text = "FoodDiaryStorage"
            break
        }
    }
}
